# Update recalculated NATMI TPM statistics for Cd28-Cd80 LR-pair sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.225147333333333
$ws.Range("H2").Value = 3.675442
$ws.Range("I2").Value = 0.2944933560673559
$ws.Range("J2").Value = 0.2944933560673559
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.9620323333333333
$ws.Range("N2").Value = 2.886097
$ws.Range("O2").Value = 0.007300799806346596
$ws.Range("P2").Value = 0.007300799806346596
$ws.Range("Q2").Value = 1.178631347763778
$ws.Range("R2").Value = 10.607682129874
$ws.Range("S2").Value = 0.002150037036946911
$ws.Range("T2").Value = 0.002150037036946912
$ws.Range("G3").Value = 1.225147333333333
$ws.Range("H3").Value = 3.675442
$ws.Range("I3").Value = 0.2944933560673559
$ws.Range("J3").Value = 0.2944933560673559
$ws.Range("O3").Value = 0.008497341958187089
$ws.Range("P3").Value = 0.008497341958187089
$ws.Range("Q3").Value = 1.371799511045555
$ws.Range("R3").Value = 12.34619559941
$ws.Range("S3").Value = 0.002502410750918474
$ws.Range("T3").Value = 0.002502410750918474
$ws.Range("G4").Value = 1.225147333333333
$ws.Range("H4").Value = 3.675442
$ws.Range("I4").Value = 0.2944933560673559
$ws.Range("J4").Value = 0.2944933560673559
$ws.Range("M4").Value = 30.69779366666667
$ws.Range("N4").Value = 92.093381
$ws.Range("O4").Value = 0.2329635276189966
$ws.Range("P4").Value = 0.2329635276189966
$ws.Range("Q4").Value = 37.60932004993356
$ws.Range("R4").Value = 338.483880449402
$ws.Range("S4").Value = 0.06860621108980848
$ws.Range("T4").Value = 0.06860621108980848
$ws.Range("G5").Value = 1.225147333333333
$ws.Range("H5").Value = 3.675442
$ws.Range("I5").Value = 0.2944933560673559
$ws.Range("J5").Value = 0.2944933560673559
$ws.Range("M5").Value = 0.4642756666666667
$ws.Range("N5").Value = 1.392827
$ws.Range("O5").Value = 0.003523357354889427
$ws.Range("P5").Value = 0.003523357354889427
$ws.Range("Q5").Value = 0.5688060949482223
$ws.Range("R5").Value = 5.119254854534001
$ws.Range("S5").Value = 0.00103760533206599
$ws.Range("T5").Value = 0.00103760533206599
$ws.Range("G6").Value = 1.225147333333333
$ws.Range("H6").Value = 3.675442
$ws.Range("I6").Value = 0.2944933560673559
$ws.Range("J6").Value = 0.2944933560673559
$ws.Range("M6").Value = 87.34624366666667
$ws.Range("N6").Value = 262.038731
$ws.Range("O6").Value = 0.6628648713262609
$ws.Range("P6").Value = 0.6628648713262609
$ws.Range("Q6").Value = 107.0120175049002
$ws.Range("R6").Value = 963.108157544102
$ws.Range("S6").Value = 0.1952093005760266
$ws.Range("T6").Value = 0.1952093005760266
$ws.Range("G7").Value = 1.225147333333333
$ws.Range("H7").Value = 3.675442
$ws.Range("I7").Value = 0.2944933560673559
$ws.Range("J7").Value = 0.2944933560673559
$ws.Range("M7").Value = 11.180767
$ws.Range("N7").Value = 33.54230099999999
$ws.Range("O7").Value = 0.08485010193531929
$ws.Range("P7").Value = 0.08485010193531929
$ws.Range("Q7").Value = 13.69808687467133
$ws.Range("R7").Value = 123.282781872042
$ws.Range("S7").Value = 0.02498779128158943
$ws.Range("T7").Value = 0.02498779128158943
$ws.Range("G8").Value = 0.005333666666666667
$ws.Range("I8").Value = 0.001282073881300198
$ws.Range("J8").Value = 0.001282073881300198
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.9620323333333333
$ws.Range("N8").Value = 2.886097
$ws.Range("O8").Value = 0.007300799806346596
$ws.Range("P8").Value = 0.007300799806346596
$ws.Range("Q8").Value = 0.005131159788555556
$ws.Range("R8").Value = 0.046180438097
$ws.Range("S8").Value = 0.000009360164744318514
$ws.Range("T8").Value = 0.000009360164744318516
$ws.Range("G9").Value = 0.005333666666666667
$ws.Range("I9").Value = 0.001282073881300198
$ws.Range("J9").Value = 0.001282073881300198
$ws.Range("O9").Value = 0.008497341958187089
$ws.Range("P9").Value = 0.008497341958187089
$ws.Range("Q9").Value = 0.005972115456111112
$ws.Range("R9").Value = 0.05374903910500001
$ws.Range("S9").Value = 0.00001089422018506795
$ws.Range("T9").Value = 0.00001089422018506795
$ws.Range("G10").Value = 0.005333666666666667
$ws.Range("I10").Value = 0.001282073881300198
$ws.Range("J10").Value = 0.001282073881300198
$ws.Range("M10").Value = 30.69779366666667
$ws.Range("N10").Value = 92.093381
$ws.Range("O10").Value = 0.2329635276189966
$ws.Range("P10").Value = 0.2329635276189966
$ws.Range("Q10").Value = 0.1637317988201111
$ws.Range("R10").Value = 1.473586189381
$ws.Range("S10").Value = 0.0002986764540558729
$ws.Range("T10").Value = 0.0002986764540558729
$ws.Range("G11").Value = 0.005333666666666667
$ws.Range("I11").Value = 0.001282073881300198
$ws.Range("J11").Value = 0.001282073881300198
$ws.Range("M11").Value = 0.4642756666666667
$ws.Range("N11").Value = 1.392827
$ws.Range("O11").Value = 0.003523357354889427
$ws.Range("P11").Value = 0.003523357354889427
$ws.Range("Q11").Value = 0.002476291647444445
$ws.Range("R11").Value = 0.022286624827
$ws.Range("S11").Value = 0.000004517204439190688
$ws.Range("T11").Value = 0.000004517204439190688
$ws.Range("G12").Value = 0.005333666666666667
$ws.Range("I12").Value = 0.001282073881300198
$ws.Range("J12").Value = 0.001282073881300198
$ws.Range("M12").Value = 87.34624366666667
$ws.Range("N12").Value = 262.038731
$ws.Range("O12").Value = 0.6628648713262609
$ws.Range("P12").Value = 0.6628648713262609
$ws.Range("Q12").Value = 0.4658757483034445
$ws.Range("R12").Value = 4.192881734731
$ws.Range("S12").Value = 0.0008498417383588158
$ws.Range("T12").Value = 0.0008498417383588158
$ws.Range("G13").Value = 0.005333666666666667
$ws.Range("I13").Value = 0.001282073881300198
$ws.Range("J13").Value = 0.001282073881300198
$ws.Range("M13").Value = 11.180767
$ws.Range("N13").Value = 33.54230099999999
$ws.Range("O13").Value = 0.08485010193531929
$ws.Range("P13").Value = 0.08485010193531929
$ws.Range("Q13").Value = 0.05963448425566666
$ws.Range("R13").Value = 0.536710358301
$ws.Range("S13").Value = 0.0001087840995169322
$ws.Range("T13").Value = 0.0001087840995169322
$ws.Range("G14").Value = 0.4766303333333334
$ws.Range("H14").Value = 1.429891
$ws.Range("I14").Value = 0.1145694584217375
$ws.Range("J14").Value = 0.1145694584217375
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.9620323333333333
$ws.Range("N14").Value = 2.886097
$ws.Range("O14").Value = 0.007300799806346596
$ws.Range("P14").Value = 0.007300799806346596
$ws.Range("Q14").Value = 0.4585337917141112
$ws.Range("R14").Value = 4.126804125427
$ws.Range("S14").Value = 0.0008364486798586554
$ws.Range("T14").Value = 0.0008364486798586556
$ws.Range("G15").Value = 0.4766303333333334
$ws.Range("H15").Value = 1.429891
$ws.Range("I15").Value = 0.1145694584217375
$ws.Range("J15").Value = 0.1145694584217375
$ws.Range("O15").Value = 0.008497341958187089
$ws.Range("P15").Value = 0.008497341958187089
$ws.Range("Q15").Value = 0.5336837786172223
$ws.Range("R15").Value = 4.803154007555001
$ws.Range("S15").Value = 0.0009735358661738012
$ws.Range("T15").Value = 0.0009735358661738013
$ws.Range("G16").Value = 0.4766303333333334
$ws.Range("H16").Value = 1.429891
$ws.Range("I16").Value = 0.1145694584217375
$ws.Range("J16").Value = 0.1145694584217375
$ws.Range("M16").Value = 30.69779366666667
$ws.Range("N16").Value = 92.093381
$ws.Range("O16").Value = 0.2329635276189966
$ws.Range("P16").Value = 0.2329635276189966
$ws.Range("Q16").Value = 14.63149962794122
$ws.Range("R16").Value = 131.683496651471
$ws.Range("S16").Value = 0.02669050519132593
$ws.Range("T16").Value = 0.02669050519132593
$ws.Range("G17").Value = 0.4766303333333334
$ws.Range("H17").Value = 1.429891
$ws.Range("I17").Value = 0.1145694584217375
$ws.Range("J17").Value = 0.1145694584217375
$ws.Range("M17").Value = 0.4642756666666667
$ws.Range("N17").Value = 1.392827
$ws.Range("O17").Value = 0.003523357354889427
$ws.Range("P17").Value = 0.003523357354889427
$ws.Range("Q17").Value = 0.2212878657618889
$ws.Range("R17").Value = 1.991590791857
$ws.Range("S17").Value = 0.0004036691439759272
$ws.Range("T17").Value = 0.0004036691439759273
$ws.Range("G18").Value = 0.4766303333333334
$ws.Range("H18").Value = 1.429891
$ws.Range("I18").Value = 0.1145694584217375
$ws.Range("J18").Value = 0.1145694584217375
$ws.Range("M18").Value = 87.34624366666667
$ws.Range("N18").Value = 262.038731
$ws.Range("O18").Value = 0.6628648713262609
$ws.Range("P18").Value = 0.6628648713262609
$ws.Range("Q18").Value = 41.6318692342579
$ws.Range("R18").Value = 374.6868231083211
$ws.Range("S18").Value = 0.07594406931464442
$ws.Range("T18").Value = 0.07594406931464444
$ws.Range("G19").Value = 0.4766303333333334
$ws.Range("H19").Value = 1.429891
$ws.Range("I19").Value = 0.1145694584217375
$ws.Range("J19").Value = 0.1145694584217375
$ws.Range("M19").Value = 11.180767
$ws.Range("N19").Value = 33.54230099999999
$ws.Range("O19").Value = 0.08485010193531929
$ws.Range("P19").Value = 0.08485010193531929
$ws.Range("Q19").Value = 5.329092702132334
$ws.Range("R19").Value = 47.961834319191
$ws.Range("S19").Value = 0.009721230225758751
$ws.Range("T19").Value = 0.009721230225758753
$ws.Range("G20").Value = 2.453075333333333
$ws.Range("H20").Value = 7.359226
$ws.Range("I20").Value = 0.5896551116296064
$ws.Range("J20").Value = 0.5896551116296064
$ws.Range("K20").Value = 3
$ws.Range("L20").Value = 1
$ws.Range("M20").Value = 0.9620323333333333
$ws.Range("N20").Value = 2.886097
$ws.Range("O20").Value = 0.007300799806346596
$ws.Range("P20").Value = 0.007300799806346596
$ws.Range("Q20").Value = 2.359937786769111
$ws.Range("R20").Value = 21.239440080922
$ws.Range("S20").Value = 0.004304953924796711
$ws.Range("T20").Value = 0.004304953924796711
$ws.Range("G21").Value = 2.453075333333333
$ws.Range("H21").Value = 7.359226
$ws.Range("I21").Value = 0.5896551116296064
$ws.Range("J21").Value = 0.5896551116296064
$ws.Range("O21").Value = 0.008497341958187089
$ws.Range("P21").Value = 0.008497341958187089
$ws.Range("Q21").Value = 2.746712539192222
$ws.Range("R21").Value = 24.72041285273
$ws.Range("S21").Value = 0.005010501120909746
$ws.Range("T21").Value = 0.005010501120909746
$ws.Range("G22").Value = 2.453075333333333
$ws.Range("H22").Value = 7.359226
$ws.Range("I22").Value = 0.5896551116296064
$ws.Range("J22").Value = 0.5896551116296064
$ws.Range("M22").Value = 30.69779366666667
$ws.Range("N22").Value = 92.093381
$ws.Range("O22").Value = 0.2329635276189966
$ws.Range("P22").Value = 0.2329635276189966
$ws.Range("Q22").Value = 75.30400043145623
$ws.Range("R22").Value = 677.736003883106
$ws.Range("S22").Value = 0.1373681348838063
$ws.Range("T22").Value = 0.1373681348838063
$ws.Range("G23").Value = 2.453075333333333
$ws.Range("H23").Value = 7.359226
$ws.Range("I23").Value = 0.5896551116296064
$ws.Range("J23").Value = 0.5896551116296064
$ws.Range("M23").Value = 0.4642756666666667
$ws.Range("N23").Value = 1.392827
$ws.Range("O23").Value = 0.003523357354889427
$ws.Range("P23").Value = 0.003523357354889427
$ws.Range("Q23").Value = 1.138903185766889
$ws.Range("R23").Value = 10.250128671902
$ws.Range("S23").Value = 0.00207756567440832
$ws.Range("T23").Value = 0.00207756567440832
$ws.Range("G24").Value = 2.453075333333333
$ws.Range("H24").Value = 7.359226
$ws.Range("I24").Value = 0.5896551116296064
$ws.Range("J24").Value = 0.5896551116296064
$ws.Range("M24").Value = 87.34624366666667
$ws.Range("N24").Value = 262.038731
$ws.Range("O24").Value = 0.6628648713262609
$ws.Range("P24").Value = 0.6628648713262609
$ws.Range("Q24").Value = 214.2669157980229
$ws.Range("R24").Value = 1928.402242182206
$ws.Range("S24").Value = 0.3908616596972311
$ws.Range("T24").Value = 0.3908616596972311
$ws.Range("G25").Value = 2.453075333333333
$ws.Range("H25").Value = 7.359226
$ws.Range("I25").Value = 0.5896551116296064
$ws.Range("J25").Value = 0.5896551116296064
$ws.Range("M25").Value = 11.180767
$ws.Range("N25").Value = 33.54230099999999
$ws.Range("O25").Value = 0.08485010193531929
$ws.Range("P25").Value = 0.08485010193531929
$ws.Range("Q25").Value = 27.42726373544733
$ws.Range("R25").Value = 246.845373619026
$ws.Range("S25").Value = 0.05003229632845418
$ws.Range("T25").Value = 0.05003229632845418
